$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Each entry: row index, Coin (B), Link (C), Price (D), Volume(1h) (E)
$rows = @(
    @(2, "Bitcoin", "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc", "37.831.39", "  +0.05%  "),
    @(3, "Ethereum", "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth", "2.077.03", "  -1.00%  "),
    @(4, "TetherUSD", "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt", "1.00", "  +0.00%  "),
    @(5, "BNB", "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb", "233.10", "  +0.37%  "),
    @(6, "XRP", "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp", "0.625", "  +0.19%  "),
    @(7, "Solana", "https://coinranking.com/coin/zNZHO_Sjf+solana-sol", "59.26", "  +2.09%  "),
    @(8, "USDC", "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc", "1.00", "  -0.03%  "),
    @(9, "Cardano", "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada", "0.394", "  +1.50%  "),
    @(10, "Dogecoin", "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge", "0.0789", "  +1.27%  "),
    @(11, "TRON", "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx", "0.107", "  +1.83%  "),
    @(12, "WrappedliquidstakedEther2.0", "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth", "2.383.21", "  -0.68%  "),
    @(13, "Chainlink", "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link", "14.78", "  +1.92%  "),
    @(14, "Avalanche", "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax", "21.24", "  +0.61%  "),
    @(15, "Polygon", "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic", "0.773", "  +0.66%  "),
    @(16, "Polkadot", "https://coinranking.com/coin/25W7FG7om+polkadot-dot", "5.35", "  +2.25%  "),
    @(17, "WrappedEther", "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth", "2.072.87", "  -1.18%  "),
    @(18, "WrappedBTC", "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc", "37.725.75", "  -0.01%  "),
    @(19, "Uniswap", "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni", "6.15", "  -0.16%  "),
    @(20, "Litecoin", "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc", "71.59", "  +1.33%  "),
    @(21, "ShibaInu", "https://coinranking.com/coin/xz24e0BjL+shibainu-shib", "0.0₃0852", "  +3.62%  "),
    @(22, "BitcoinCash", "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch", "228.31", "  +0.12%  "),
    @(23, "Dai", "https://coinranking.com/coin/MoTuySvg7+dai-dai", "0.999", "  -0.02%  "),
    @(24, "PancakeSwap", "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake", "2.42", "  +1.58%  "),
    @(25, "Toncoin", "https://coinranking.com/coin/67YlI0K1b+toncoin-ton", "2.37", "  -1.10%  "),
    @(26, "Monero", "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr", "170.79", "  +1.75%  "),
    @(27, "Cosmos", "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom", "9.20", "  +2.79%  "),
    @(28, "Kaspa", "https://coinranking.com/coin/V8GxkwWow+kaspa-kas", "0.134", "  -4.54%  "),
    @(29, "ImmutableX", "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx", "1.42", "  -0.82%  "),
    @(30, "EthereumClassic", "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc", "19.49", "  +0.22%  "),
    @(31, "Stellar", "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm", "0.121", "  +1.57%  "),
    @(32, "Filecoin", "https://coinranking.com/coin/ymQub4fuB+filecoin-fil", "4.73", "  +2.32%  "),
    @(33, "InternetComputer(DFINITY)", "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp", "4.75", "  +3.69%  "),
    @(34, "Hedera", "https://coinranking.com/coin/jad286TjB+hedera-hbar", "0.0633", "  +1.08%  "),
    @(35, "LidoDAOToken", "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo", "2.50", "  -0.13%  "),
    @(36, "WEMIXToken", "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix", "1.82", "  -0.20%  "),
    @(37, "RenderToken", "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr", "3.40", "  -0.52%  "),
    @(38, "BinanceUSD", "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd", "1.00", "  +0.23%  "),
    @(39, "THORChain", "https://coinranking.com/coin/ybmU-kKU+thorchain-rune", "5.43", "  +0.29%  "),
    @(40, "Cronos", "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro", "0.0985", "  -0.90%  "),
    @(41, "Aave", "https://coinranking.com/coin/ixgUfzmLR+aave-aave", "99.19", "  +1.20%  "),
    @(42, "VeChain", "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet", "0.0216", "  +0.88%  "),
    @(43, "HuobiToken", "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht", "2.88", "  -1.95%  "),
    @(44, "InjectiveProtocol", "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj", "16.70", "  +6.62%  "),
    @(45, "Maker", "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr", "1.444.62", "  -0.91%  "),
    @(46, "TrustWalletToken", "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt", "1.16", "  -0.26%  "),
    @(47, "FTXToken", "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt", "4.23", "  +4.71%  "),
    @(48, "ARBITRUM", "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb", "1.06", "  +0.38%  "),
    @(49, "FraxShare", "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs", "7.39", "  +0.45%  "),
    @(50, "MXToken", "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx", "3.01", "  -0.34%  "),
    @(51, "RocketPoolETH", "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth", "2.268.19", "  -0.96%  ")
)

foreach ($entry in $rows) {
    $r = $entry[0]
    $ws.Cells.Item($r, 2).Value = $entry[1]
    $ws.Cells.Item($r, 3).Value = $entry[2]

    # Price (D) and Volume (E) values must stay plain text, matching the
    # original inlineStr cells. Force text format so values that look
    # numeric (e.g. "1.00", "233.10") aren't auto-converted to numbers,
    # then restore the default "Normal" style so no stray style index is
    # left on the cell (keeps the XML identical to the original layout).
    $dCell = $ws.Cells.Item($r, 4)
    $dCell.NumberFormat = "@"
    $dCell.Value = $entry[3]
    $dCell.Style = "Normal"

    $eCell = $ws.Cells.Item($r, 5)
    $eCell.NumberFormat = "@"
    $eCell.Value = $entry[4]
    $eCell.Style = "Normal"
}
